# Updates cryptos list values (Price / Volume(1h) columns, and a row swap
# between Kaspa and RenzoRestakedETH) to match the latest scrape.
#
# Note: several "Price" values are plain numeric-looking strings (e.g.
# "675.61"). The source workbook stores them as *text* (inline strings),
# so when writing them back through the Excel object model we prefix the
# value with a leading apostrophe to force text entry and avoid Excel
# silently converting them to floating point numbers (which would also
# mangle trailing zeros, e.g. "8.20" -> 8.2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-TextValue 'D2' '69.289.23'
$ws.Range('E2').Value = '  -0.08%  '

# Row 3 - Ethereum
Set-TextValue 'D3' '3.672.68'
$ws.Range('E3').Value = '  -0.32%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.05%  '

# Row 5 - BNB
Set-TextValue 'D5' '675.61'
$ws.Range('E5').Value = '  -1.05%  '

# Row 6 - Solana
Set-TextValue 'D6' '158.44'
$ws.Range('E6').Value = '  -2.41%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  -0.10%  '

# Row 8 - XRP
$ws.Range('E8').Value = '  -1.35%  '

# Row 9 - Dogecoin
$ws.Range('E9').Value = '  -1.64%  '

# Row 10 - Toncoin
$ws.Range('E10').Value = '  -5.51%  '

# Row 11 - Cardano
$ws.Range('E11').Value = '  -2.42%  '

# Row 12 - ShibaInu
Set-TextValue 'D12' '0.0000232'
$ws.Range('E12').Value = '  -3.65%  '

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue 'D13' '4.292.33'
$ws.Range('E13').Value = '  -0.33%  '

# Row 14 - Avalanche
Set-TextValue 'D14' '32.38'
$ws.Range('E14').Value = '  -3.88%  '

# Row 15 - WrappedEther
Set-TextValue 'D15' '3.680.80'
$ws.Range('E15').Value = '  -0.05%  '

# Row 16 - WrappedBTC
Set-TextValue 'D16' '69.219.08'

# Row 17 - TRON
$ws.Range('E17').Value = '  +1.54%  '

# Row 18 - Chainlink
$ws.Range('E18').Value = '  -1.89%  '

# Row 19 - Polkadot
$ws.Range('E19').Value = '  -2.82%  '

# Row 20 - BitcoinCash
Set-TextValue 'D20' '467.28'
$ws.Range('E20').Value = '  -3.06%  '

# Row 21 - Uniswap
$ws.Range('E21').Value = '  +0.78%  '

# Row 22 - Polygon
$ws.Range('E22').Value = '  -2.80%  '

# Row 23 - Litecoin
$ws.Range('E23').Value = '  -0.59%  '

# Row 24 - WrappedeETH
Set-TextValue 'D24' '3.818.96'
$ws.Range('E24').Value = '  -0.29%  '

# Row 25 - Dai
$ws.Range('E25').Value = '  -0.05%  '

# Row 26 - PEPE
$ws.Range('E26').Value = '  -6.59%  '

# Row 27 - InternetComputer(DFINITY)
Set-TextValue 'D27' '10.88'
$ws.Range('E27').Value = '  -5.53%  '

# Row 28 - RenderToken
$ws.Range('E28').Value = '  -4.70%  '

# Row 30 - Fetch.AI
$ws.Range('E30').Value = '  -4.71%  '

# Row 31 - NEARProtocol
Set-TextValue 'D31' '6.62'
$ws.Range('E31').Value = '  -3.21%  '

# Row 32 - Binance-PegBSC-USD
$ws.Range('E32').Value = '  -0.13%  '

# Row 33 - ImmutableX
Set-TextValue 'D33' '1.98'
$ws.Range('E33').Value = '  -4.81%  '

# Row 34 - EthereumClassic
Set-TextValue 'D34' '26.89'
$ws.Range('E34').Value = '  -0.81%  '

# Row 35/36 - Kaspa and RenzoRestakedETH swapped order
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D35' '0.162'
$ws.Range('E35').Value = '  -4.08%  '

$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue 'D36' '3.666.00'
$ws.Range('E36').Value = '  +0.30%  '

# Row 37 - Aptos
Set-TextValue 'D37' '8.20'
$ws.Range('E37').Value = '  -3.50%  '

# Row 38 - Filecoin
Set-TextValue 'D38' '6.25'
$ws.Range('E38').Value = '  -1.03%  '

# Row 40 - FirstDigitalUSD
$ws.Range('E40').Value = '  -0.03%  '

# Row 41 - Stacks
$ws.Range('E41').Value = '  -4.34%  '

# Row 42 - Hedera
Set-TextValue 'D42' '0.0902'
$ws.Range('E42').Value = '  -3.76%  '

# Row 43 - Monero
Set-TextValue 'D43' '173.41'
$ws.Range('E43').Value = '  +6.81%  '

# Row 44 - Mantle
Set-TextValue 'D44' '0.941'

# Row 45 - OKB
Set-TextValue 'D45' '47.66'
$ws.Range('E45').Value = '  -1.33%  '

# Row 46 - InjectiveProtocol
Set-TextValue 'D46' '28.26'
$ws.Range('E46').Value = '  -6.08%  '

# Row 47 - FLOKI
Set-TextValue 'D47' '0.000277'
$ws.Range('E47').Value = '  -3.28%  '

# Row 48 - dogwifhat
$ws.Range('E48').Value = '  -5.03%  '

# Row 49 - ONDO
$ws.Range('E49').Value = '  -4.46%  '

# Row 50 - SuiNetwork
$ws.Range('E50').Value = '  -4.42%  '

# Row 51 - Cosmos
$ws.Range('E51').Value = '  -3.07%  '
